# ---------------------------------------------------------------------------
# feat: add 2022-Q3 data
#
# 1. "总计" (summary) sheet: insert a new top data row for "2022-Q3"
#    (count=8, value=3.98) and shift the existing 2022-Q2 / 2022-Q1 / 2021-Q4
#    rows down by one.
# 2. Insert a brand new "2022-Q3" worksheet (duplicated from "2022-Q2" so it
#    keeps identical layout/formatting) positioned right before "2022-Q2",
#    then overwrite it with the Q3 fund holdings (8 funds).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1 - "总计" sheet: shift rows down (bottom-up) and write the new row.
# ---------------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")

# row4 (2021-Q4) -> row5
$totalWs.Cells.Item(4,1).Copy($totalWs.Cells.Item(5,1))
$totalWs.Cells.Item(4,2).Copy($totalWs.Cells.Item(5,2))
$totalWs.Cells.Item(4,3).Copy($totalWs.Cells.Item(5,3))
$totalWs.Cells.Item(4,4).Copy($totalWs.Cells.Item(5,4))
$totalWs.Cells.Item(5,1).Value = 3

# row3 (2022-Q1) -> row4
$totalWs.Cells.Item(3,1).Copy($totalWs.Cells.Item(4,1))
$totalWs.Cells.Item(3,2).Copy($totalWs.Cells.Item(4,2))
$totalWs.Cells.Item(3,3).Copy($totalWs.Cells.Item(4,3))
$totalWs.Cells.Item(3,4).Copy($totalWs.Cells.Item(4,4))
$totalWs.Cells.Item(4,1).Value = 2

# row2 (2022-Q2) -> row3
$totalWs.Cells.Item(2,1).Copy($totalWs.Cells.Item(3,1))
$totalWs.Cells.Item(2,2).Copy($totalWs.Cells.Item(3,2))
$totalWs.Cells.Item(2,3).Copy($totalWs.Cells.Item(3,3))
$totalWs.Cells.Item(2,4).Copy($totalWs.Cells.Item(3,4))
$totalWs.Cells.Item(3,1).Value = 1

# row2 becomes the new 2022-Q3 summary row
$totalWs.Cells.Item(2,1).Value = 0
$totalWs.Cells.Item(2,2).Value = "2022-Q3"
$totalWs.Cells.Item(2,3).Value = 8
$totalWs.Cells.Item(2,4).Value = 3.98

# ---------------------------------------------------------------------------
# Step 2 - add the "2022-Q3" detail sheet right before "2022-Q2", duplicating
# the "2022-Q2" sheet so formatting/dimensions match, then overwrite values.
# ---------------------------------------------------------------------------
$q2Ws = $wb.Worksheets.Item("2022-Q2")
$q2Index = $q2Ws.Index
$q2Ws.Copy($q2Ws)
# The duplicate is inserted immediately before the source sheet, so it now
# sits at the source's former index (use position rather than the
# auto-generated "<name> (2)" title to avoid relying on that naming scheme).
$q3Ws = $wb.Worksheets.Item($q2Index)
$q3Ws.Name = "2022-Q3"

# "2022-Q2" currently has 4 rows of data (rows 2-4); Q3 needs 8 (rows 2-9).
# Extend the sheet by replicating the last data row's formatting downward.
for ($r = 5; $r -le 9; $r++) {
    $q3Ws.Range("A4:H4").Copy($q3Ws.Range("A" + $r + ":H" + $r))
}

# Row 2
$q3Ws.Cells.Item(2,1).Value = 0
$q3Ws.Cells.Item(2,2).Value = "'007455"
$q3Ws.Cells.Item(2,3).Value = "富国蓝筹精选股票（QDII）人民币"
$q3Ws.Cells.Item(2,4).Value = "'13.62"
$q3Ws.Cells.Item(2,5).Value = "'85.59"
$q3Ws.Cells.Item(2,6).Value = "'5.62"
$q3Ws.Cells.Item(2,7).Value = "'0.7654"
$q3Ws.Cells.Item(2,8).Value = 3

# Row 3
$q3Ws.Cells.Item(3,1).Value = 1
$q3Ws.Cells.Item(3,2).Value = "'010583"
$q3Ws.Cells.Item(3,3).Value = "富国蓝筹精选股票（QDII）美元"
$q3Ws.Cells.Item(3,4).Value = "'13.62"
$q3Ws.Cells.Item(3,5).Value = "'85.59"
$q3Ws.Cells.Item(3,6).Value = "'5.62"
$q3Ws.Cells.Item(3,7).Value = "'0.7654"
$q3Ws.Cells.Item(3,8).Value = 3

# Row 4
$q3Ws.Cells.Item(4,1).Value = 2
$q3Ws.Cells.Item(4,2).Value = "'000934"
$q3Ws.Cells.Item(4,3).Value = "国富大中华精选混合（QDII）"
$q3Ws.Cells.Item(4,4).Value = "'19.83"
$q3Ws.Cells.Item(4,5).Value = "'72.45"
$q3Ws.Cells.Item(4,6).Value = "'3.31"
$q3Ws.Cells.Item(4,7).Value = "'0.6564"
$q3Ws.Cells.Item(4,8).Value = 2

# Row 5
$q3Ws.Cells.Item(5,1).Value = 3
$q3Ws.Cells.Item(5,2).Value = "'006370"
$q3Ws.Cells.Item(5,3).Value = "国富大中华精选混合（QDII）美元"
$q3Ws.Cells.Item(5,4).Value = "'19.83"
$q3Ws.Cells.Item(5,5).Value = "'72.45"
$q3Ws.Cells.Item(5,6).Value = "'3.31"
$q3Ws.Cells.Item(5,7).Value = "'0.6564"
$q3Ws.Cells.Item(5,8).Value = 2

# Row 6
$q3Ws.Cells.Item(6,1).Value = 4
$q3Ws.Cells.Item(6,2).Value = "'010671"
$q3Ws.Cells.Item(6,3).Value = "景顺长城大中华混合（QDII）美元A"
$q3Ws.Cells.Item(6,4).Value = "'9.42"
$q3Ws.Cells.Item(6,5).Value = "'70.56"
$q3Ws.Cells.Item(6,6).Value = "'4.25"
$q3Ws.Cells.Item(6,7).Value = "'0.4004"
$q3Ws.Cells.Item(6,8).Value = 7

# Row 7
$q3Ws.Cells.Item(7,1).Value = 5
$q3Ws.Cells.Item(7,2).Value = "'262001"
$q3Ws.Cells.Item(7,3).Value = "景顺长城大中华混合（QDII）人民币A"
$q3Ws.Cells.Item(7,4).Value = "'9.42"
$q3Ws.Cells.Item(7,5).Value = "'70.56"
$q3Ws.Cells.Item(7,6).Value = "'4.25"
$q3Ws.Cells.Item(7,7).Value = "'0.4004"
$q3Ws.Cells.Item(7,8).Value = 7

# Row 8
$q3Ws.Cells.Item(8,1).Value = 6
$q3Ws.Cells.Item(8,2).Value = "'100055"
$q3Ws.Cells.Item(8,3).Value = "富国全球科技互联网股票（QDII）"
$q3Ws.Cells.Item(8,4).Value = "'3.95"
$q3Ws.Cells.Item(8,5).Value = "'86.97"
$q3Ws.Cells.Item(8,6).Value = "'5.03"
$q3Ws.Cells.Item(8,7).Value = "'0.1987"
$q3Ws.Cells.Item(8,8).Value = 3

# Row 9
$q3Ws.Cells.Item(9,1).Value = 7
$q3Ws.Cells.Item(9,2).Value = "'457001"
$q3Ws.Cells.Item(9,3).Value = "国富亚洲机会股票（QDII）"
$q3Ws.Cells.Item(9,4).Value = "'3.80"
$q3Ws.Cells.Item(9,5).Value = "'83.80"
$q3Ws.Cells.Item(9,6).Value = "'3.48"
$q3Ws.Cells.Item(9,7).Value = "'0.1322"
$q3Ws.Cells.Item(9,8).Value = 2
